# Generate Report for Handback
# Updates the localization-status report for the file
# "c8e39951-57c8-4bf4-8150-3dbc85837702.md" now that its handback has been
# processed: status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the Latest Handback DateTime stamps are
# refreshed, and the (now stale) Error Detail text is cleared.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns for the c8e39951 row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

# --- zh-cn sheet: Status / Latest Handback DateTime / Error Detail ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("K3").Value = "2016-08-21 02:52:27"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: Status / Latest Handback DateTime / Error Detail ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("K3").Value = "2016-08-21 02:52:34"
$wsDeDe.Range("P3").Value = ""
